# Re-apply the canonical per-match row ordering for "Germany Bundesliga I".
# The underlying data source re-synced with corrected match ids, which shuffled
# which match (columns B:AD) appears on which existing spreadsheet row while the
# running index in column A (and the row position itself) stays put. We therefore
# snapshot each affected rows B:AD payload first (so sources are read before any
# destination is overwritten), then write every row to its corrected match data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 77, 78, 79 exchange their match data (id in col B is the giveaway).
$snap77 = $ws.Range("B77:AD77").Value2
$snap78 = $ws.Range("B78:AD78").Value2
$snap79 = $ws.Range("B79:AD79").Value2
$ws.Range("B77:AD77").Value2 = $snap78
$ws.Range("B78:AD78").Value2 = $snap79
$ws.Range("B79:AD79").Value2 = $snap77

# Rows 85, 88 exchange their match data (id in col B is the giveaway).
$snap85 = $ws.Range("B85:AD85").Value2
$snap88 = $ws.Range("B88:AD88").Value2
$ws.Range("B85:AD85").Value2 = $snap88
$ws.Range("B88:AD88").Value2 = $snap85

# Rows 93, 96, 95 exchange their match data (id in col B is the giveaway).
$snap93 = $ws.Range("B93:AD93").Value2
$snap96 = $ws.Range("B96:AD96").Value2
$snap95 = $ws.Range("B95:AD95").Value2
$ws.Range("B93:AD93").Value2 = $snap96
$ws.Range("B96:AD96").Value2 = $snap95
$ws.Range("B95:AD95").Value2 = $snap93

# Rows 137, 138 exchange their match data (id in col B is the giveaway).
$snap137 = $ws.Range("B137:AD137").Value2
$snap138 = $ws.Range("B138:AD138").Value2
$ws.Range("B137:AD137").Value2 = $snap138
$ws.Range("B138:AD138").Value2 = $snap137

# Rows 140, 144 exchange their match data (id in col B is the giveaway).
$snap140 = $ws.Range("B140:AD140").Value2
$snap144 = $ws.Range("B144:AD144").Value2
$ws.Range("B140:AD140").Value2 = $snap144
$ws.Range("B144:AD144").Value2 = $snap140

# Rows 141, 143 exchange their match data (id in col B is the giveaway).
$snap141 = $ws.Range("B141:AD141").Value2
$snap143 = $ws.Range("B143:AD143").Value2
$ws.Range("B141:AD141").Value2 = $snap143
$ws.Range("B143:AD143").Value2 = $snap141

# Rows 154, 155, 158, 156 exchange their match data (id in col B is the giveaway).
$snap154 = $ws.Range("B154:AD154").Value2
$snap155 = $ws.Range("B155:AD155").Value2
$snap158 = $ws.Range("B158:AD158").Value2
$snap156 = $ws.Range("B156:AD156").Value2
$ws.Range("B154:AD154").Value2 = $snap155
$ws.Range("B155:AD155").Value2 = $snap158
$ws.Range("B158:AD158").Value2 = $snap156
$ws.Range("B156:AD156").Value2 = $snap154

# Rows 173, 176 exchange their match data (id in col B is the giveaway).
$snap173 = $ws.Range("B173:AD173").Value2
$snap176 = $ws.Range("B176:AD176").Value2
$ws.Range("B173:AD173").Value2 = $snap176
$ws.Range("B176:AD176").Value2 = $snap173

# Rows 273, 274, 275 exchange their match data (id in col B is the giveaway).
$snap273 = $ws.Range("B273:AD273").Value2
$snap274 = $ws.Range("B274:AD274").Value2
$snap275 = $ws.Range("B275:AD275").Value2
$ws.Range("B273:AD273").Value2 = $snap274
$ws.Range("B274:AD274").Value2 = $snap275
$ws.Range("B275:AD275").Value2 = $snap273

